$wb = $excel.ActiveWorkbook

function Set-RowData {
    param(
        $ws,
        [int]$row,
        [string]$timeVal,
        [string]$colB,
        [string]$colC,
        [string]$colD,
        [string]$colE,
        [int]$colF,
        [string]$colG,
        [int]$colH,
        [int]$colI
    )

    # Column A: timestamp, copy the date/time number format from the row above
    # so the new cell matches the existing column styling exactly.
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
    $ws.Cells.Item($row, 1).Value2 = $timeVal

    # Columns B-E: text payload fields.
    $ws.Cells.Item($row, 2).Value2 = $colB
    $ws.Cells.Item($row, 3).Value2 = $colC
    $ws.Cells.Item($row, 4).Value2 = $colD
    $ws.Cells.Item($row, 5).Value2 = $colE

    # Column F: plain integer.
    $ws.Cells.Item($row, 6).Value2 = $colF

    # Column G: very large decimal number - assign via Value2 using the exact
    # decimal text so the stored double matches precisely, then strip any
    # auto-applied "scientific" number format so the cell stays unstyled.
    $ws.Cells.Item($row, 7).Value2 = $colG
    $ws.Cells.Item($row, 7).ClearFormats()

    # Columns H and I: plain integers.
    $ws.Cells.Item($row, 8).Value2 = $colH
    $ws.Cells.Item($row, 9).Value2 = $colI
}

$ws1 = $wb.Worksheets.Item("ROW35-FE-LIFTER")
Set-RowData $ws1 53 "45750.8550134375" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x72" "0xd" 400 "5.68631262647114E+23" 370 13

$ws2 = $wb.Worksheets.Item("ROW35-MID-LIFTER")
Set-RowData $ws2 53 "45750.70503052083" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x6e" "0xe" 400 "5.68631262647114E+23" 366 14

$ws3 = $wb.Worksheets.Item("ROW02-FE-LIFTER")
Set-RowData $ws3 53 "45750.84803054398" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x01,0x72" "0x3" 400 "5.68631262647114E+23" 370 3

$ws4 = $wb.Worksheets.Item("ROW02-MID-LIFTER")
Set-RowData $ws4 53 "45750.90638702546" "0x01,0x90" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x01,0x6e" "0x3" 400 "9.85046333984776E+23" 366 3
